$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row: "<name>_old" -> "<name>_FV2310", "<name>_new" -> "<name>_FV2404" ---
$oldCols = @("A","B","C","D","E","F","G","H","I","J")
$newCols = @("L","M","N","O","P","Q","R","S","T","U")

foreach ($col in $oldCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value()) -replace "_old$", "_FV2310"
}

foreach ($col in $newCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value()) -replace "_new$", "_FV2404"
}

# --- 2. Turn the used range into an Excel Table (ListObject) ---
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U55"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# --- 3. Freeze the header row (pane split below row 1) ---
$null = $ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
